$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new rows 54-61 (columns A, B, C) in the precise order that
# reproduces the shared-string table order from the target workbook.
$ws.Range("A55").Value = "more images"
$ws.Range("B55").Value = "weitere Bilder"
$ws.Range("C55").Value = "mehr Bilder"
$ws.Range("B56").Value = "Reihenfolge festzulegen"
$ws.Range("C56").Value = "Reihenfolge zu arrangieren"
$ws.Range("A56").Value = "arrange their order"
$ws.Range("A57").Value = "start"
$ws.Range("A58").Value = "information"
$ws.Range("B57").Value = "Starten"
$ws.Range("B58").Value = "Informationen"
$ws.Range("B59").Value = "Fokus"
$ws.Range("C57").Value = "anfangen"
$ws.Range("C58").Value = "informieren"
$ws.Range("C59").Value = "fokussieren"
$ws.Range("A54").Value = "Start & End Frames"
$ws.Range("B54").Value = "Anfangs- und Endframe"
$ws.Range("A60").Value = "you'd like to keep"
$ws.Range("B60").Value = "die Sie...behalten möchten"
$ws.Range("C60").Value = "die Sie behalten möchten"
$ws.Range("A61").Value = "We recommend downloading the ones"
$ws.Range("B61").Value = "Wir empfehlen diejenigen herunterzuladen"
$ws.Range("C61").Value = "Wir empfehlen, diejenigen…herunterzuladen"
$ws.Range("C54").Value = "Anfangs- & End-Frames"
$ws.Range("A59").Value = "focus"

# Apply the same style (vertical-center + wrap text) used throughout the
# table to the new rows.
$ws.Range("A54:C61").WrapText = $true
$ws.Range("A54:C61").VerticalAlignment = -4108

# Row 61 has an explicit custom row height (ht="30" in target XML).
$ws.Rows("61").RowHeight = 30

# Column D gets the same styling applied (but no values) for rows 53-59,
# matching the source formatting extension seen in the target workbook.
$ws.Range("D53:D59").WrapText = $true
$ws.Range("D53:D59").VerticalAlignment = -4108

# Update page setup to portrait orientation.
$ws.PageSetup.Orientation = 1

# Update the view's top-left cell and active selection/cell.
$excel.ActiveWindow.ScrollRow = 28
$ws.Range("H51").Select() | Out-Null
